$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G56").Value = 3.7
$ws.Range("H56").Value = 3.5
$ws.Range("K56").Value = 2.38
$ws.Range("M56").Value = 1.03
$ws.Range("N56").Value = 15
$ws.Range("O56").Value = 1.18
$ws.Range("P56").Value = 4.5
$ws.Range("Q56").Value = 1.65
$ws.Range("R56").Value = 2.2
$ws.Range("S56").Value = 1.3
$ws.Range("T56").Value = 3.4
$ws.Range("U56").Value = 1.57
$ws.Range("V56").Value = 2.25
$ws.Range("W56").Value = 15
$ws.Range("X56").Value = 23
$ws.Range("AC56").Value = 15
$ws.Range("AD56").Value = 7.5
$ws.Range("AH56").Value = 9.5
$ws.Range("AI56").Value = 10
$ws.Range("AL56").Value = 13
$ws.Range("AM56").Value = 21
$ws.Range("AP56").Value = 23
$ws.Range("AR56").Value = 67
$ws.Range("AT56").Value = 3.4
$ws.Range("AV56").Value = 41
$ws.Range("AY56").Value = 17
$ws.Range("BA56").Value = 41
$ws.Range("BB56").Value = 101
$ws.Range("J61").Value = 2.8
$ws.Range("M61").Value = 1.1
$ws.Range("N61").Value = 6.5
$ws.Range("O61").Value = 1.4
$ws.Range("P61").Value = 2.52
$ws.Range("Q61").Value = 2.15
$ws.Range("R61").Value = 1.55
$ws.Range("T61").Value = 2.35
$ws.Range("U61").Value = 1.85
$ws.Range("V61").Value = 1.75
$ws.Range("W61").Value = 6.5
$ws.Range("X61").Value = 10
$ws.Range("Z61").Value = 22
$ws.Range("AA61").Value = 19.5
$ws.Range("AB61").Value = 32
$ws.Range("AC61").Value = 7.2
$ws.Range("AE61").Value = 15.5
$ws.Range("AG61").Value = 800
$ws.Range("AL61").Value = 37
$ws.Range("AN61").Value = 3.95
$ws.Range("AP61").Value = 21
$ws.Range("AQ61").Value = 45
$ws.Range("AR61").Value = 80
$ws.Range("AT61").Value = 2.32
$ws.Range("AU61").Value = 7.1
$ws.Range("AV61").Value = 70
$ws.Range("AW61").Value = 5.2
$ws.Range("AY61").Value = 28
